# Replace curly/smart quotes with straight single quotes in five English
# story lines (column C) of the sheet, per the 20210731 "画中人" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "No one would offer the Infected any help. We were no volunteer army. We had nothing like the 'Brave Wok' from the age of Mikhail the Huge Ears. We had no city of our own, and we hardly even had fields to hill up." + "`n"
$ws.Range("C4").Value = "''The Originium grows on our skin, and the scant few weapons we hold are stark. The snow in our mouths melts into our water, our stomachs are filled with grass seed and bark. '" + "`n"
$ws.Range("C50").Value = "[name=`"Talulah`"] " + ([char]0x2014) + "'Aegis.'  'Aegis' are attacking the sentries here." + "`n"
$ws.Range("C82").Value = "[name=`"Infected Fighter`"] A 'decision' means nothing to us. We came back with you all the way from the Northwest... and it wasn" + ([char]0x2019) + "t for any 'decision.'" + "`n"
$ws.Range("C90").Value = "[name=`"Talulah`"] It" + ([char]0x2019) + "s okay. Let me go. Less casualties this way, too.  Besides, I" + ([char]0x2019) + "m" + ([char]0x2014) + "*chuckle* 'the Deathless.'" + "`n"
